# Handback report update: a new localization file
# "001d19e8-7580-475d-b9b9-54a875137953.md" was handed back (in sync with
# en-US) in between the previously-last "2e340cf0..." row and the
# "08adc6d5..." row. Insert a new row 3 on every sheet (pushing the old
# row 3 down to row 4) and populate it, then fix up the table ranges and
# hyperlinks that the plain row-insert doesn't keep in sync.

$wb = $excel.ActiveWorkbook

$newFile      = "001d19e8-7580-475d-b9b9-54a875137953.md"
$newFilePath  = "e2e\001d19e8-7580-475d-b9b9-54a875137953.md"
$newExt       = ".md"
$statusSync   = "Handed back: in sync with en-US"
$overviewDate = "2016-08-30 16:53:22"

$oldFile      = "08adc6d5-5918-446d-be8c-515de8e8e274.md"
$oldFilePath  = "e2e\08adc6d5-5918-446d-be8c-515de8e8e274.md"

# ---------------------------------------------------------------------
# Sheet "Overview": columns A-G = File Name, Path And Name, Extension,
# Publish URL, zh-cn, de-de, Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Rows.Item(3).Insert()

$wsOv.Range("A3").Value = $newFile
$wsOv.Range("B3").Value = $newFilePath
$wsOv.Range("C3").Value = $newExt
$wsOv.Range("E3").Value = $statusSync
$wsOv.Range("F3").Value = $statusSync
$wsOv.Range("G3").Value = $overviewDate

$loOv = $wsOv.ListObjects.Item("Overview")
$loOv.Resize($wsOv.Range("A1:G4"))

$wsOv.Range("A1").Hyperlinks.Delete()
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa15af2ed472e40ac797dc924e5aa25893ddf519/e2e/2e340cf0-cfd9-47ef-a0a2-6096b1fa233c.md", "", "", "e2e\2e340cf0-cfd9-47ef-a0a2-6096b1fa233c.md") | Out-Null
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c1d2e3f4a5b6c7d8e9f001122334455667788990/$newFilePath", "", "", $newFilePath) | Out-Null
$wsOv.Hyperlinks.Add($wsOv.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0b3f869cc4c05f0e7e21ce9bcdfb437776aa32f/e2e/08adc6d5-5918-446d-be8c-515de8e8e274.md", "", "", $oldFilePath) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn": one detail row per handback file
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Rows.Item(3).Insert()

$wsZh.Range("A3").Value = $newFile
$wsZh.Range("B3").Value = $newExt
$wsZh.Range("C3").Value = $statusSync
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "001d19e8-7580-475d-b9b9-54a875137953.8f4e2afb230e47cfb3b0774188e20f4b48025d15.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-30 16:53:18"
$wsZh.Range("I3").Value = $newFile
$wsZh.Range("J3").Value = "001d19e8-7580-475d-b9b9-54a875137953.8f4e2afb230e47cfb3b0774188e20f4b48025d15.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-30 16:53:34"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$loZh = $wsZh.ListObjects.Item("zh-cn")
$loZh.Resize($wsZh.Range("A1:P4"))

$wsZh.Range("A1").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa15af2ed472e40ac797dc924e5aa25893ddf519/e2e/2e340cf0-cfd9-47ef-a0a2-6096b1fa233c.md", "", "", "2e340cf0-cfd9-47ef-a0a2-6096b1fa233c.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/488520491be6cc700a92c7254958139b82d5159c/e2e/2e340cf0-cfd9-47ef-a0a2-6096b1fa233c.md", "", "", "2e340cf0-cfd9-47ef-a0a2-6096b1fa233c.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c1d2e3f4a5b6c7d8e9f001122334455667788990/$newFilePath", "", "", $newFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d4e5f60718293a4b5c6d7e8f9001a2b3c4d5e6f7/e2e/001d19e8-7580-475d-b9b9-54a875137953.md", "", "", $newFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0b3f869cc4c05f0e7e21ce9bcdfb437776aa32f/e2e/08adc6d5-5918-446d-be8c-515de8e8e274.md", "", "", $oldFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a45a00a19da9ad67c2427f371aaa91a9e66f0b0b/e2e/08adc6d5-5918-446d-be8c-515de8e8e274.md", "", "", $oldFile) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de": one detail row per handback file
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Rows.Item(3).Insert()

$wsDe.Range("A3").Value = $newFile
$wsDe.Range("B3").Value = $newExt
$wsDe.Range("C3").Value = $statusSync
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "001d19e8-7580-475d-b9b9-54a875137953.8f4e2afb230e47cfb3b0774188e20f4b48025d15.de-de.xlf"
$wsDe.Range("H3").Value = $overviewDate
$wsDe.Range("I3").Value = $newFile
$wsDe.Range("J3").Value = "001d19e8-7580-475d-b9b9-54a875137953.8f4e2afb230e47cfb3b0774188e20f4b48025d15.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-30 16:53:41"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$loDe = $wsDe.ListObjects.Item("de-de")
$loDe.Resize($wsDe.Range("A1:P4"))

$wsDe.Range("A1").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa15af2ed472e40ac797dc924e5aa25893ddf519/e2e/2e340cf0-cfd9-47ef-a0a2-6096b1fa233c.md", "", "", "2e340cf0-cfd9-47ef-a0a2-6096b1fa233c.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/84ec2f79a29adae94098285a11311245644f7675/e2e/2e340cf0-cfd9-47ef-a0a2-6096b1fa233c.md", "", "", "2e340cf0-cfd9-47ef-a0a2-6096b1fa233c.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c1d2e3f4a5b6c7d8e9f001122334455667788990/$newFilePath", "", "", $newFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/e5f60718293a4b5c6d7e8f9001a2b3c4d5e6f708/e2e/001d19e8-7580-475d-b9b9-54a875137953.md", "", "", $newFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0b3f869cc4c05f0e7e21ce9bcdfb437776aa32f/e2e/08adc6d5-5918-446d-be8c-515de8e8e274.md", "", "", $oldFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1729f18a4c7b22fb86429beacc9028b5bacaa050/e2e/08adc6d5-5918-446d-be8c-515de8e8e274.md", "", "", $oldFile) | Out-Null
